# Fruta / hortaliza, semanal
# Insert a new weekly record at row 27 (pushing the existing rows 27-41 down
# to 28-42), then append three new weekly records as rows 43-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 27 (Tuna / Primera, week of 2022-01-06) ------------
$ws.Rows.Item(27).Insert()

$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D27").Value = 44567
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = 100112027
$ws.Range("G27").Value = "Melón"
$ws.Range("H27").Value = "Tuna"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 9000
$ws.Range("M27").Value = 8500
$ws.Range("N27").Value = "$/caja 18 unidades"
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("P27").Value = 472
$ws.Range("Q27").Value = 18
$ws.Range("R27").Value = "Hortaliza"

# --- Append new row 43 (Calameño / Primera, week of 2022-01-07) --------
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C43").Value = "Arica y Parinacota"
$ws.Range("D43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D43").Value = 44568
$ws.Range("E43").Value = 15
$ws.Range("F43").Value = 100112027
$ws.Range("G43").Value = "Melón"
$ws.Range("H43").Value = "Calameño"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 70
$ws.Range("K43").Value = 9000
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = 9500
$ws.Range("N43").Value = "$/caja 18 unidades"
$ws.Range("O43").Value = "Región de Arica y Parinacota"
$ws.Range("P43").Value = 528
$ws.Range("Q43").Value = 18
$ws.Range("R43").Value = "Hortaliza"

# --- Append new row 44 (Calameño / Segunda, week of 2022-01-07) -------
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C44").Value = "Arica y Parinacota"
$ws.Range("D44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D44").Value = 44568
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = 100112027
$ws.Range("G44").Value = "Melón"
$ws.Range("H44").Value = "Calameño"
$ws.Range("I44").Value = "Segunda"
$ws.Range("J44").Value = 100
$ws.Range("K44").Value = 5000
$ws.Range("L44").Value = 5500
$ws.Range("M44").Value = 5250
$ws.Range("N44").Value = "$/caja 24 unidades"
$ws.Range("O44").Value = "Región de Arica y Parinacota"
$ws.Range("P44").Value = 219
$ws.Range("Q44").Value = 24
$ws.Range("R44").Value = "Hortaliza"

# --- Append new row 45 (Tuna / Primera, week of 2022-01-07) -----------
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C45").Value = "Arica y Parinacota"
$ws.Range("D45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D45").Value = 44568
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 100112027
$ws.Range("G45").Value = "Melón"
$ws.Range("H45").Value = "Tuna"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 80
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 11000
$ws.Range("M45").Value = 10500
$ws.Range("N45").Value = "$/caja 18 unidades"
$ws.Range("O45").Value = "Región de Arica y Parinacota"
$ws.Range("P45").Value = 583
$ws.Range("Q45").Value = 18
$ws.Range("R45").Value = "Hortaliza"
